# Applies the "Physics -> Chemistry" content edit to the document.
#
# Strategy:
#  - Plain text swaps (title, author, email, body sentences) are done by
#    locating the old text with Range.Find.Execute (search-only, no
#    replacement argument) and then assigning the new string straight to
#    Range.Text. Assigning .Text (rather than passing the replacement
#    through Find.Execute's Replace arg) avoids Word's AutoCorrect
#    "smart quotes" from mangling the straight apostrophes that are in
#    the target copy.
#  - The very last body paragraph (the "Summary" paragraph) also needs a
#    <w:lastRenderedPageBreak/> moved to a new spot mid-paragraph and one
#    of its sentences dropped outright, so it is rebuilt wholesale with
#    Range.InsertXML using a minimal WordprocessingML package fragment.

$d = $word.ActiveDocument

function Replace-DocText($find, $replace) {
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($find)
    if (-not $found) {
        throw "Text not found: $find"
    }
    $rng.Text = $replace
}

# --- Title ---------------------------------------------------------------
Replace-DocText "The Fabric of the Universe-  A Physicist's Perspective" `
    "The Enchanting World of Chemistry: Unveiling Nature's Symphony of Elements"

# --- Byline ----------------------------------------------------------------
Replace-DocText "Isaac Newton" "Dr. Emily Carter"

# --- Email -----------------------------------------------------------------
Replace-DocText "isaac.newton@physics.com" "ecarter@highschoolchem.edu"

# --- Body paragraph (sentence by sentence, formatting/breaks untouched) ---
Replace-DocText `
    "The cosmos, a symphony of celestial bodies, captivates our imagination and fuels scientific inquiry" `
    "In the vast tapestry of sciences, chemistry stands out as a symphony of elements, a harmonious dance of molecules, and a vibrant narrative of matter"

Replace-DocText `
    " Physics unlocks the secrets of the universe, delving into the fundamental laws that govern matter and energy, space and time" `
    " It delves into the intricate relationships between substances, unveiling the secrets of their composition, properties, and interactions"

Replace-DocText `
    " One of the most prominent theories in physics, quantum mechanics, has profoundly shaped our understanding of the universe at its smallest scales" `
    " Chemistry is the language of the natural world, spoken in the patterns of atoms, the eloquence of chemical bonds, and the ever-changing states of matter"

Replace-DocText `
    " This fascinating realm, ruled by enigmatic particles and forces, has unveiled a hidden tapestry of quantum superposition, entanglement, and wave-particle duality, forever altering our perception of reality" `
    " Through chemistry, we gain the power to decipher nature's enigmatic whispers, unravel the complexities of our world, and harness its boundless potential for innovation"

Replace-DocText `
    "As we peer deeper into the vastness of the cosmos, a symphony of celestial bodies reveals the profound unity underlying the universe's diverse phenomena" `
    "With each new discovery, chemistry opens doors to uncharted territories of knowledge, revealing the hidden wonders of the universe"

Replace-DocText `
    " Einstein's theory of general relativity, a captivating tapestry of spacetime curvature, gravity, and the cosmic dance of celestial objects, has transformed our understanding of gravity and the cosmos" `
    " It empowers us to decode the intricate workings of life, unlock the secrets of disease, and devise ingenious solutions to global challenges"

Replace-DocText `
    " General relativity's elegance and predictive power have enabled us to unlock the mysteries of black holes, gravitational waves, and the expansion of the universe, offering a glimpse into the breathtaking vastness of existence" `
    " Chemistry is the key to unraveling the mysteries of the cosmos, understanding the intricate mechanisms of our bodies, and developing revolutionary technologies that shape our future"

Replace-DocText `
    "Delving into the realm of subatomic particles, the Standard Model of Physics gracefully orchestrates the intricate ballet of fundamental forces and particles, providing a comprehensive framework that encompasses the electromagnetic, weak, and strong interactions" `
    "As we delve deeper into the enchanting world of chemistry, we embark on a journey of exploration, experimentation, and enlightenment"

Replace-DocText `
    " This symphony of subatomic interactions forms the foundation of matter and energy, dictating the properties of atoms, molecules, and the world we experience" `
    " We become alchemists, transforming ordinary substances into extraordinary materials, unlocking the secrets of chemical reactions, and witnessing the magic of transformations"

Replace-DocText `
    " From the birth of stars to the fusion within them and the radiant melodies of atomic transitions, physics unlocks the secrets of energy transformation and the symphony of the universe" `
    " Chemistry invites us to explore the boundless possibilities of matter, to create new substances, and to understand the interconnectedness of all things"

# --- Summary paragraph: rebuilt wholesale -----------------------------
# New text (with the lastRenderedPageBreak now sitting mid-sentence,
# right before "enchanting world of chemistry, we embark...") and one
# whole sentence ("Thus, we embark on an endless quest...") dropped.
$summaryParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$summaryRange = $summaryParagraph.Range

$summaryXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>The realm of chemistry is an enchanting tapestry of elements, molecules, and interactions, offering a symphony of knowledge about the composition, properties, and transformations of matter</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> It empowers us to understand the natural world, decode the complexities of life, and devise innovative solutions to global challenges</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> As we delve deeper into the </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t>enchanting world of chemistry, we embark on a journey of exploration, experimentation, and enlightenment, unlocking the secrets of matter and shaping our future</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$summaryRange.InsertXML($summaryXml)
